$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-19 10:52:02"
$wsZhCn.Range("H2").Value = "2016-03-19 10:52:20"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-19 10:52:05"
$wsDeDe.Range("H2").Value = "2016-03-19 10:52:25"
